# Portfolio data sheet: add the first saved asset row (row 2) under the
# existing header row (row 1). Mirrors a single "ETF prova" ETF position
# persisted by the app (camelCase JSON fields -> columns A..T).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1              # Id
$ws.Range("B2").Value = "ETF"          # category
$ws.Range("C2").Value = "ETF prova"    # assetName
$ws.Range("D2").Value = "Fineco"       # position

$ws.Range("E2").Value = 1              # riskLevel

# ticker / isin were not provided for this asset. Writing a lone "'" forces
# Excel to treat the cell as (empty) text rather than dropping it, then the
# style is reset to the sheet default so no quote-prefix formatting lingers
# on the cell - this reproduces the blank-but-present text cells the app
# originally wrote for these two fields.
$ws.Range("F2").Value = "'"            # ticker
$ws.Range("F2").Style = "Normal"
$ws.Range("G2").Value = "'"            # isin
$ws.Range("G2").Style = "Normal"

$ws.Range("H2").Value = "22/05/2000"   # createdAt (kept as literal text)
$ws.Range("I2").Value = 2              # createdAmount
$ws.Range("J2").Value = 22000          # createdUnitPrice
$ws.Range("K2").Value = 44000          # createdTotalValue

$ws.Range("L2").Value = "22/7/2024"    # updatedAt (kept as literal text)
$ws.Range("M2").Value = 3              # updatedAmount
$ws.Range("N2").Value = 27000          # updatedUnitPrice
$ws.Range("O2").Value = 81000          # updatedTotalValue

$ws.Range("P2").Value = "'"            # accumulationPlan (blank text)
$ws.Range("P2").Style = "Normal"

$ws.Range("Q2").Value = 0              # accumulationAmount
$ws.Range("R2").Value = 0              # incomePerYear
$ws.Range("S2").Value = 0              # rentalIncome
$ws.Range("T2").Value = "Prova asset"  # note
